$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.660.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.596.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0842"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.821.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.590.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.634.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0736"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.14%  "
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "208.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("E21").Value = "  +4.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.28%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("E28").Value = "  -0.90%  "
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("E30").Value = "  +1.77%  "
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.291.04"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("E35").Value = "  -5.40%  "
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("E39").Value = "  -1.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +16.73%  "
$ws.Range("E41").Value = "  +1.01%  "
$ws.Range("E42").Value = "  -0.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.781"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.732.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.77%  "
$ws.Range("E48").Value = "  +1.30%  "
$ws.Range("E49").Value = "  +0.92%  "
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.78%  "
